# Append a new departure row (row 7) to the "Main Data" sheet, mirroring
# the layout of the existing rows (columns A..M, with K and L/M used for
# the blank "difference" separator + value as in rows 2-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Sunday, Jan 08"
$ws.Range("C7").Value = "9:35 PM"
$ws.Range("D7").Value = "FR6640"
$ws.Range("E7").Value = "London"
$ws.Range("F7").Value = "(LTN)"
$ws.Range("G7").Value = "Ryanair "
$ws.Range("H7").Value = "B738"
$ws.Range("I7").Value = "(EI-EMH)"
$ws.Range("J7").Value = "10:21 PM"
$ws.Range("L7").Value = "0 hours, 46 minutes"
